$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'27.356.38"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +2.35%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'1.661.94"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +1.30%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  -0.37%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'220.09"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +1.05%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'0.507"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  +0.64%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'  -0.37%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = "'  +1.22%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'0.0627"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  +0.26%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'19.97"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  +4.39%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('E11').Value = "'  +1.02%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'1.894.66"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  +1.35%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'1.666.28"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  +1.71%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('E14').Value = "'  +1.24%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'0.533"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +1.17%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'67.24"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  +3.96%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'27.337.88"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E18').Value = "'  +0.52%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'222.25"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  +3.25%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('E20').Value = "'  -0.37%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('E21').Value = "'  +9.15%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'4.45"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  +1.91%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'2.51"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  +6.06%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'9.29"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  +0.15%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'147.20"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  +1.26%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = "'  -0.27%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'7.44"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  +3.75%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Value = "'  +1.05%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'16.05"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  +2.63%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('E30').Value = "'  +1.36%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D32').Value = "'3.39"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  +0.05%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'3.01"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  +0.14%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('E34').Value = "'  +2.41%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'1.264.80"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  -1.83%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').Value = "'  +0.61%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Value = "'  +0.03%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'0.537"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  -0.16%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'0.835"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  +2.28%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('E40').Value = "'  -0.37%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'0.815"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  +1.32%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('E42').Value = "'  +2.66%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'1.805.98"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  +1.50%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('E44').Value = "'  -4.03%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'61.83"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  +1.37%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'92.56"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  +0.88%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').Value = "'  +1.47%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = "'  -0.52%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'0.0985"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  +2.10%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = "'  +0.48%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('E51').Value = "'  +0.42%  "
$ws.Range('E51').Style = 'Normal'
